$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.839.42'
$ws.Range("E2").Value = '  +0.44%  '

$ws.Range("D3").Value = '2.304.55'
$ws.Range("E3").Value = '  +1.06%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = "'301.67"
$ws.Range("E5").Value = '  -1.19%  '

$ws.Range("D6").Value = "'96.19"
$ws.Range("E6").Value = '  -0.26%  '

$ws.Range("E7").Value = '  +0.54%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").Value = "'0.495"
$ws.Range("E9").Value = '  -1.03%  '

$ws.Range("D10").Value = "'34.66"
$ws.Range("E10").Value = '  -2.34%  '

$ws.Range("D11").Value = "'19.25"
$ws.Range("E11").Value = '  +5.38%  '

$ws.Range("D12").Value = "'0.0790"
$ws.Range("E12").Value = '  -0.09%  '

$ws.Range("D13").Value = "'0.119"
$ws.Range("E13").Value = '  +0.26%  '

$ws.Range("E14").Value = '  +0.89%  '

$ws.Range("D15").Value = '2.667.56'
$ws.Range("E15").Value = '  +1.23%  '

$ws.Range("D16").Value = '2.312.57'
$ws.Range("E16").Value = '  +1.91%  '

$ws.Range("E17").Value = '  +0.81%  '

$ws.Range("D18").Value = '42.756.22'
$ws.Range("E18").Value = '  +0.41%  '

$ws.Range("D19").Value = "'12.20"
$ws.Range("E19").Value = '  -5.64%  '

$ws.Range("E20").Value = '  -0.50%  '

$ws.Range("D21").Value = "'6.01"
$ws.Range("E21").Value = '  +0.14%  '

$ws.Range("E22").Value = '  +1.00%  '

$ws.Range("D23").Value = "'2.27"
$ws.Range("E23").Value = '  +7.40%  '

$ws.Range("D24").Value = "'235.28"
$ws.Range("E24").Value = '  -0.32%  '

$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("D26").Value = "'2.41"
$ws.Range("E26").Value = '  -2.07%  '

$ws.Range("E27").Value = '  -2.80%  '

$ws.Range("E28").Value = '  +14.81%  '

$ws.Range("D29").Value = "'165.61"
$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("D30").Value = "'9.08"
$ws.Range("E30").Value = '  +0.50%  '

$ws.Range("E31").Value = '  -2.80%  '

$ws.Range("E32").Value = '  -0.04%  '

$ws.Range("E33").Value = '  +0.60%  '

$ws.Range("D34").Value = "'17.65"
$ws.Range("E34").Value = '  +0.60%  '

$ws.Range("E35").Value = '  -6.28%  '

$ws.Range("D36").Value = "'0.0701"
$ws.Range("E36").Value = '  +1.81%  '

$ws.Range("E38").Value = '  -0.52%  '

$ws.Range("E39").Value = '  +0.87%  '

$ws.Range("E40").Value = '  -0.68%  '

$ws.Range("E41").Value = '  +0.44%  '

$ws.Range("D42").Value = "'20.35"
$ws.Range("E42").Value = '  +12.74%  '

$ws.Range("D43").Value = '1.969.75'
$ws.Range("E43").Value = '  -1.50%  '

$ws.Range("D44").Value = "'10.49"
$ws.Range("E44").Value = '  +5.32%  '

$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("E46").Value = '  -1.87%  '

$ws.Range("E47").Value = '  -0.09%  '

$ws.Range("B48").Value = 'HuobiToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D48").Value = "'2.84"
$ws.Range("E48").Value = '  -0.81%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.531.40'
$ws.Range("E49").Value = '  +1.09%  '

$ws.Range("D50").Value = "'53.33"
$ws.Range("E50").Value = '  -0.20%  '

$ws.Range("D51").Value = "'71.52"
$ws.Range("E51").Value = '  +0.35%  '
